$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Make ProductLoanInput the active sheet/tab
$ws1.Activate()

# Update row 6: "Currency" -> "currency", "US Dollar " -> "US Dollar"
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"

# Row B6 picks up a fresh default-font style with a light-green fill,
# matching the highlighted "currency" row header styling used for A6.
$ws1.Range("B6").Style = "Normal"
$ws1.Range("B6").Interior.Color = 5296274

# Update the visible selection on the ProductLoanInput sheet to the edited row
$ws1.Range("A6:B6").Select()

$wb.Save()
